$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 9: add PriceChange (X9) and UpDown (Y9)
$ws.Range("X9").Value = -1.7200020000000222
$ws.Range("Y9").Value = "Down"

# Add new row 10
$ws.Range("A10").Value = 42653.880208333336
$ws.Range("B10").Value = -8
$ws.Range("C10").Value = "Sell"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = "Random"
$ws.Range("Q10").Value = 44.409433632991338
$ws.Range("R10").Value = 1.8
$ws.Range("S10").Value = 0.0926
$ws.Range("T10").Value = -0.0094
$ws.Range("U10").Value = 5.87
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = 0

$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("S9").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("T9").Copy()
$ws.Range("T10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
